$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) values for the refreshed
# cryptocurrency snapshot. Column D values look like plain numbers, so the
# NumberFormat is forced to text before assignment to preserve exact
# formatting (thousand-dot separators, trailing zeros, subscript digits),
# then the style is reset back to Normal so no stray formatting is left
# behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.675.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.998.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "379.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.546"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.87"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.139"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.64%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0846"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.475.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.998.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.966"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.727.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "264.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +16.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.172"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "26.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.105"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.96%  "

$ws.Range("E35").Value = "  -3.93%  "

$ws.Range("E36").Value = "  +4.63%  "

$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.50%  "

$ws.Range("E42").Value = "  +1.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.285"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +19.75%  "

$ws.Range("E46").Value = "  -2.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.95%  "

$ws.Range("E48").Value = "  +1.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.051.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.02%  "

$ws.Range("E50").Value = "  +10.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.880"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.15%  "

Write-Host "Updated cryptos list"
